$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simplify the login-flow description text: drop the "관리자 및 회원"/"관리자"
# (admin-specific) wording so the steps read as generic login steps. These three
# pieces of text each appear three times in the sheet (대여소 등록 / 대여소 리스트 조회 /
# 대여소 상세 정보 조회 blocks), so every occurrence is updated individually.
$text2 = "2. 로그인, 회원 가입 메뉴를 띄운다."
$text3 = "3. 로그인 메뉴를 누른다."
$text4 = "4. 로그인 화면을 띄운다."

$ws.Range("E13").Value = $text2
$ws.Range("H13").Value = $text2
$ws.Range("E21").Value = $text2
$ws.Range("H21").Value = $text2
$ws.Range("E30").Value = $text2
$ws.Range("H30").Value = $text2

$ws.Range("D14").Value = $text3
$ws.Range("G14").Value = $text3
$ws.Range("D22").Value = $text3
$ws.Range("G22").Value = $text3
$ws.Range("D31").Value = $text3
$ws.Range("G31").Value = $text3

$ws.Range("E14").Value = $text4
$ws.Range("H14").Value = $text4
$ws.Range("E22").Value = $text4
$ws.Range("H22").Value = $text4
$ws.Range("E31").Value = $text4
$ws.Range("H31").Value = $text4

# --- Remove stray leftover formatting in columns A:C (rows 35:36) and the
# trailing blank formatted row 49 so the sheet's used range shrinks back down
# to C2:H48.
$ws.Range("A35:C36").Clear()
$ws.Range("D49:E49").Clear()

# --- Update the saved view/selection state.
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("G26:H26").Select()
